$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.335781415968434
$ws.Range("D2").Value = 4.548792062674039
$ws.Range("E2").Value = 10.86996394917689
$ws.Range("F2").Value = 25.28316119936956
$ws.Range("G2").Value = 31.00171573934332
$ws.Range("H2").Value = 14.42545230254646
$ws.Range("I2").Value = 20.83873490581467
$ws.Range("K2").Value = 12.86289807999148
$ws.Range("M2").Value = 15.20128943422548
$ws.Range("N2").Value = 17.22642167257143
$ws.Range("C3").Value = 3.263275551446421
$ws.Range("D3").Value = 4.544437965874929
$ws.Range("E3").Value = 10.7694390102926
$ws.Range("F3").Value = 25.17309191483605
$ws.Range("G3").Value = 30.74402645771804
$ws.Range("H3").Value = 14.45218063625729
$ws.Range("I3").Value = 20.80593083347798
$ws.Range("K3").Value = 12.37888381533752
$ws.Range("M3").Value = 14.91959767325879
$ws.Range("N3").Value = 17.29171022561286
$ws.Range("C4").Value = 3.217400953408053
$ws.Range("D4").Value = 4.541847298415709
$ws.Range("E4").Value = 10.71136834086205
$ws.Range("F4").Value = 25.11433420896835
$ws.Range("G4").Value = 30.59835237976664
$ws.Range("H4").Value = 14.47212863619439
$ws.Range("I4").Value = 20.79247663658521
$ws.Range("K4").Value = 12.07409173371841
$ws.Range("M4").Value = 14.74781653073327
$ws.Range("N4").Value = 17.33368307631869
$ws.Range("C5").Value = 3.198380299835382
$ws.Range("D5").Value = 4.540812848657552
$ws.Range("E5").Value = 10.68864568291932
$ws.Range("F5").Value = 25.0926255951351
$ws.Range("G5").Value = 30.54220516487712
$ws.Range("H5").Value = 14.48114397800461
$ws.Range("I5").Value = 20.78867648444703
$ws.Range("K5").Value = 11.94817103844263
$ws.Range("M5").Value = 14.67820914260216
$ws.Range("N5").Value = 17.35126296523614
$ws.Range("C6").Value = 3.195202663231361
$ws.Range("D6").Value = 4.540642376977014
$ws.Range("E6").Value = 10.68493010862904
$ws.Range("F6").Value = 25.08915632942038
$ws.Range("G6").Value = 30.53307783036331
$ws.Range("H6").Value = 14.48269441154562
$ws.Range("I6").Value = 20.78814708351215
$ws.Range("K6").Value = 11.92716454582114
$ws.Range("M6").Value = 14.66667762257767
$ws.Range("N6").Value = 17.35421086265879
$ws.Range("C7").Value = 3.217145735926675
$ws.Range("D7").Value = 4.541833260764585
$ws.Range("E7").Value = 10.71105805357463
$ws.Range("F7").Value = 25.11403236819426
$ws.Range("G7").Value = 30.5975820623976
$ws.Range("H7").Value = 14.47224663552445
$ws.Range("I7").Value = 20.79241857411201
$ws.Range("K7").Value = 12.07240019389706
$ws.Range("M7").Value = 14.74687605207403
$ws.Range("N7").Value = 17.33391823713713
$ws.Range("C8").Value = 3.311071211906547
$ws.Range("D8").Value = 4.547273864108568
$ws.Range("E8").Value = 10.83455962508573
$ws.Range("F8").Value = 25.24338862424596
$ws.Range("G8").Value = 30.91029900998304
$ws.Range("H8").Value = 14.43393248608948
$ws.Range("I8").Value = 20.82603711581553
$ws.Range("K8").Value = 12.69770044239603
$ws.Range("M8").Value = 15.10397626216311
$ws.Range("N8").Value = 17.2485427570645
$ws.Range("C9").Value = 3.483909213785595
$ws.Range("D9").Value = 4.558580115060735
$ws.Range("E9").Value = 11.1045928405611
$ws.Range("F9").Value = 25.56618085831171
$ws.Range("G9").Value = 31.62004648398505
$ws.Range("H9").Value = 14.38698263543629
$ws.Range("I9").Value = 20.94490571897563
$ws.Range("K9").Value = 13.85591892234691
$ws.Range("M9").Value = 15.80933475919606
$ws.Range("N9").Value = 17.09601226559922
$ws.Range("C10").Value = 3.603244700487151
$ws.Range("D10").Value = 4.567250007162525
$ws.Range("E10").Value = 11.31831409423794
$ws.Range("F10").Value = 25.84402541306034
$ws.Range("G10").Value = 32.19575265906731
$ws.Range("H10").Value = 14.36982148329478
$ws.Range("I10").Value = 21.06423297996626
$ws.Range("K10").Value = 14.65650917560252
$ws.Range("M10").Value = 16.32509247963839
$ws.Range("N10").Value = 16.9929279767273
$ws.Range("C11").Value = 3.655736139465261
$ws.Range("D11").Value = 4.571267713949116
$ws.Range("E11").Value = 11.41848694881348
$ws.Range("F11").Value = 25.97887921249669
$ws.Range("G11").Value = 32.46834083011465
$ws.Range("H11").Value = 14.36580435410914
$ws.Range("I11").Value = 21.12536668808135
$ws.Range("K11").Value = 15.00825843670934
$ws.Range("M11").Value = 16.55806750180158
$ws.Range("N11").Value = 16.94796144858577
$ws.Range("C12").Value = 3.675344902967386
$ws.Range("D12").Value = 4.572799256467228
$ws.Range("E12").Value = 11.45681031490524
$ws.Range("F12").Value = 26.03112487171942
$ws.Range("G12").Value = 32.57300268579012
$ws.Range("H12").Value = 14.36482969584759
$ws.Range("I12").Value = 21.1494903515848
$ws.Range("K12").Value = 15.13956151023745
$ws.Range("M12").Value = 16.645964235053
$ws.Range("N12").Value = 16.93120937306588
$ws.Range("C13").Value = 3.671133909380323
$ws.Range("D13").Value = 4.572468968698872
$ws.Range("E13").Value = 11.44853989408501
$ws.Range("F13").Value = 26.0198209995015
$ws.Range("G13").Value = 32.55039951767937
$ws.Range("H13").Value = 14.36501527640059
$ws.Range("I13").Value = 21.14425178795083
$ws.Range("K13").Value = 15.1113689613537
$ws.Range("M13").Value = 16.62704996555806
$ws.Range("N13").Value = 16.93480498818293
$ws.Range("C14").Value = 3.657354804623516
$ws.Range("D14").Value = 4.571393514704438
$ws.Range("E14").Value = 11.42163221536566
$ws.Range("F14").Value = 25.98315407930835
$ws.Range("G14").Value = 32.47692317320571
$ws.Range("H14").Value = 14.36571320689489
$ws.Range("I14").Value = 21.12733191527671
$ws.Range("K14").Value = 15.01909939768024
$ws.Range("M14").Value = 16.5653058444366
$ws.Range("N14").Value = 16.94657772591366
$ws.Range("C15").Value = 3.648879427594206
$ws.Range("D15").Value = 4.570736072379329
$ws.Range("E15").Value = 11.40520026023328
$ws.Range("F15").Value = 25.96084697968726
$ws.Range("G15").Value = 32.43210113331583
$ws.Range("H15").Value = 14.36621192716513
$ws.Range("I15").Value = 21.11709443068191
$ws.Range("K15").Value = 14.96233167054654
$ws.Range("M15").Value = 16.52744078114986
$ws.Range("N15").Value = 16.95382474226678
$ws.Range("C16").Value = 3.599777226509742
$ws.Range("D16").Value = 4.566988883298479
$ws.Range("E16").Value = 11.31182402403161
$ws.Range("F16").Value = 25.83537953796888
$ws.Range("G16").Value = 32.17814519357668
$ws.Range("H16").Value = 14.37016041293646
$ws.Range("I16").Value = 21.06037466607921
$ws.Range("K16").Value = 14.63326163188083
$ws.Range("M16").Value = 16.30982617164686
$ws.Range("N16").Value = 16.99590530757171
$ws.Range("C17").Value = 3.569187159278201
$ws.Range("D17").Value = 4.564708647457618
$ws.Range("E17").Value = 11.25527226037266
$ws.Range("F17").Value = 25.76055079063584
$ws.Range("G17").Value = 32.02502140197485
$ws.Range("H17").Value = 14.37355452995414
$ws.Range("I17").Value = 21.02732603634645
$ws.Range("K17").Value = 14.42811974331706
$ws.Range("M17").Value = 16.1758402865994
$ws.Range("N17").Value = 17.0222129763787
$ws.Range("C18").Value = 3.551424192440975
$ws.Range("D18").Value = 4.563404084304552
$ws.Range("E18").Value = 11.22302494875999
$ws.Range("F18").Value = 25.71831007094514
$ws.Range("G18").Value = 31.93796166245298
$ws.Range("H18").Value = 14.37586335305994
$ws.Range("I18").Value = 21.0089630060879
$ws.Range("K18").Value = 14.30896178822164
$ws.Range("M18").Value = 16.09862531008915
$ws.Range("N18").Value = 17.03752591127823
$ws.Range("C19").Value = 3.545381373664545
$ws.Range("D19").Value = 4.562963593039468
$ws.Range("E19").Value = 11.21215564049368
$ws.Range("F19").Value = 25.70414637100793
$ws.Range("G19").Value = 31.90866169247342
$ws.Range("H19").Value = 14.37670627572963
$ws.Range("I19").Value = 21.002856815527
$ws.Range("K19").Value = 14.26842024147761
$ws.Range("M19").Value = 16.07245870070122
$ws.Range("N19").Value = 17.04274181474081
$ws.Range("C20").Value = 3.572461028547049
$ws.Range("D20").Value = 4.564950664724292
$ws.Range("E20").Value = 11.26126358214959
$ws.Range("F20").Value = 25.76843402673687
$ws.Range("G20").Value = 32.04121755780744
$ws.Range("H20").Value = 14.37315629914649
$ws.Range("I20").Value = 21.03077737302147
$ws.Range("K20").Value = 14.45007892373363
$ws.Range("M20").Value = 16.19011948079078
$ws.Range("N20").Value = 17.01939370941505
$ws.Range("C21").Value = 3.661409431052333
$ws.Range("D21").Value = 4.571709130646747
$ws.Range("E21").Value = 11.42952534161753
$ws.Range("F21").Value = 25.99389234283699
$ws.Range("G21").Value = 32.49846672179829
$ws.Range("H21").Value = 14.36549336369614
$ws.Range("I21").Value = 21.13227536802296
$ws.Range("K21").Value = 15.04625346752735
$ws.Range("M21").Value = 16.58345112160714
$ws.Range("N21").Value = 16.943112313587
$ws.Range("C22").Value = 3.717972441018458
$ws.Range("D22").Value = 4.576184969475476
$ws.Range("E22").Value = 11.54175155310261
$ws.Range("F22").Value = 26.14809799076635
$ws.Range("G22").Value = 32.80564267825753
$ws.Range("H22").Value = 14.36367135709289
$ws.Range("I22").Value = 21.20427879071552
$ws.Range("K22").Value = 15.42479229859025
$ws.Range("M22").Value = 16.83857863094271
$ws.Range("N22").Value = 16.89486482271954
$ws.Range("C23").Value = 3.687930596316306
$ws.Range("D23").Value = 4.573790908265737
$ws.Range("E23").Value = 11.48165919790273
$ws.Range("F23").Value = 26.06518123864444
$ws.Range("G23").Value = 32.64096790562609
$ws.Range("H23").Value = 14.36435181967733
$ws.Range("I23").Value = 21.16533479637161
$ws.Range("K23").Value = 15.2238061250323
$ws.Range("M23").Value = 16.7026178487764
$ws.Range("N23").Value = 16.92046882825236
$ws.Range("C24").Value = 3.570981460675442
$ws.Range("D24").Value = 4.564841228792172
$ws.Range("E24").Value = 11.25855407739682
$ws.Range("F24").Value = 25.76486758525682
$ws.Range("G24").Value = 32.0338922367643
$ws.Range("H24").Value = 14.37333522575374
$ws.Range("I24").Value = 21.02921503834969
$ws.Range("K24").Value = 14.44015496733709
$ws.Range("M24").Value = 16.18366442556293
$ws.Range("N24").Value = 17.02066771329001
$ws.Range("C25").Value = 3.438445800367609
$ws.Range("D25").Value = 4.555455869104498
$ws.Range("E25").Value = 11.02872410124661
$ws.Range("F25").Value = 25.47159278157201
$ws.Range("G25").Value = 31.41815736695614
$ws.Range("H25").Value = 14.3966497283175
$ws.Range("I25").Value = 20.90710230495746
$ws.Range("K25").Value = 13.55085924151522
$ws.Range("M25").Value = 15.61855781128906
$ws.Range("N25").Value = 17.13569185321555
